$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I0 and IF
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J (rows 2-10)
$values = @(7, 7, 9, 8, 8, 9, 9, 8, 5)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
